$d = $word.ActiveDocument

# 1) Professional summary paragraph: "affecting all Black and Asian-American voters" -> "affecting 50M voters"
$d.Content.Find.Execute(
    "Discovered systematic demographic coding errors affecting all Black and Asian-American voters, developed",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Discovered systematic demographic coding errors affecting 50M voters, developed",
    2) | Out-Null

# 2) Work-experience bullet: replace the phrase, then re-bold/color just the "50M" portion
#    so it matches the formatting already used for the "23%"/"64%" figures in the same run.
#    (Scoped to this paragraph only -- the professional summary above also now contains
#    "50M" and must stay a single, unformatted run.)
$bulletPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate = $d.Paragraphs.Item($i)
    if ($candidate.Range.Text -like "*Discovered systematic race coding errors affecting*") {
        $bulletPara = $candidate
        break
    }
}

$bulletRange = $bulletPara.Range.Duplicate
$bulletRange.Find.Execute(
    "affecting all Black and Asian-American voters, developed geospatial machine learning",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "affecting 50M voters, developed geospatial machine learning",
    2) | Out-Null

$numRng = $bulletPara.Range.Duplicate
$numRng.Find.Execute("50M") | Out-Null
$numRng.Bold = $true
$numRng.Font.Color = 5258796

# 3) Project impact line: "affecting all Black and Asian-American voters" -> "affecting 50M voters nationwide"
$d.Content.Find.Execute(
    "Impact: Corrected demographic data affecting all Black and Asian-American voters, improved",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Impact: Corrected demographic data affecting 50M voters nationwide, improved",
    2) | Out-Null
